$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23, shifting existing rows 23..58 down to 24..59
$ws.Rows.Item(23).Insert()

# Populate the new row 23 with its data (mirrors the template used by sibling rows)
$ws.Cells.Item(23, 1).Value = 10
$ws.Cells.Item(23, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(23, 3).Value = "La Araucanía"
$ws.Cells.Item(23, 4).Value = 44789
$ws.Cells.Item(23, 5).Value = 9
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100108
$ws.Cells.Item(23, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(23, 9).Value = 100108003
$ws.Cells.Item(23, 10).Value = "Maracuyá"
$ws.Cells.Item(23, 11).Value = "Sin especificar"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 40
$ws.Cells.Item(23, 14).Value = 36000
$ws.Cells.Item(23, 15).Value = 36000
$ws.Cells.Item(23, 16).Value = 36000
$ws.Cells.Item(23, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(23, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(23, 19).Value = 2000
$ws.Cells.Item(23, 20).Value = 18
